# Validation + Exceptions Handlers
# Adds a small "feature checklist" block (columns I:M) next to rows 9-11,
# highlights three additional self-assessment answers in green (F9, D10,
# E11), adds a width for the new column K, and updates the view (zoom /
# scroll / selection) on the first sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Highlight three more "chosen answer" cells with the existing green
#     fill used elsewhere in the grading grid (e.g. G12, D13). ---
$green = 5287936   # RGB(0,176,80) in BGR COM color order
$ws.Range("F9").Interior.Color = $green
$ws.Range("D10").Interior.Color = $green
$ws.Range("E11").Interior.Color = $green

# --- New column (K) width ---
$ws.Columns.Item(11).ColumnWidth = 12.33

# --- Feature checklist block, columns I:M ---
# Row 9 already has values; only its vertical alignment changes to "top".
$ws.Range("I9:M9").VerticalAlignment = -4160

# Row 10: existing I10:K10 get the same alignment change, and two new
# (empty) cells L10/M10 are introduced with that same formatting.
$ws.Range("I10:M10").VerticalAlignment = -4160

# Row 11: brand-new checklist entries.
$ws.Range("I11").Value = "podgląda kursów walut"
$ws.Range("J11").Value = "podgląd klientów"
$ws.Range("K11").Value = "pogdląd klienta, jego adresów oraz histori wymian"
$ws.Range("L11").Value = "histroia tranksakcji"
$ws.Range("I11:M11").VerticalAlignment = -4160

# --- Sheet view: zoom in and move the visible window / selection ---
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I14").Select()
